# Fix contact information missing from short resumes.
#
# The resume's name line ("Dheeraj Chand") was not followed by a contact
# info line. Insert a new, centered paragraph directly after it containing
# the phone / email / website / LinkedIn / location, matching the long
# resume layout.

$d = $word.ActiveDocument

$contactLine = "202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX"

# Use Find & Replace with the special "^p" paragraph-mark code so the new
# paragraph is created the way Word itself would when a user types Enter
# after the name and then types the contact line - this produces a plain
# new paragraph (centered, inheriting no bold/size run formatting) rather
# than one that clones the "Dheeraj Chand" run's direct character
# formatting.
$found = $d.Content.Find.Execute("Dheeraj Chand", $true, $false, $false, $false, $false, $true, 1, $false, "Dheeraj Chand^p" + $contactLine, 2)

if (-not $found) {
    throw "Could not find 'Dheeraj Chand' heading to insert contact info after."
}

# The newly-created paragraph should be centered, same as the name above it.
$namePara = $d.Paragraphs(1)
$contactPara = $d.Paragraphs(2)
$contactPara.Alignment = $namePara.Alignment

Write-Output "Inserted contact info paragraph after name."
